# Daily attendance processing - 2026-01-03 11:00:29
# Swap the order of "Recorded By" entries in column G:
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
